$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new text-valued cells (B:O) serialize as text, matching the existing
# inline-string cells used throughout this column range, rather than being
# auto-coerced to numbers.
$ws.Range("B35:O55").NumberFormat = "@"

# Match the bold / bordered / centered look already used by column A (id)
# for every existing data row.
$idRange = $ws.Range("A35:A55")
$idRange.Font.Bold = $true
$idRange.HorizontalAlignment = -4108
$idRange.VerticalAlignment = -4160
$idRange.Borders.LineStyle = 1
$idRange.Borders.Weight = 2

$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(35, 3).Value = "ranking"
$ws.Cells.Item(35, 4).Value = "tensorflow"
$ws.Cells.Item(35, 5).Value = "12/03/2018"
$ws.Cells.Item(35, 6).Value = "0"
$ws.Cells.Item(35, 7).Value = "0"
$ws.Cells.Item(35, 8).Value = "0"
$ws.Cells.Item(35, 9).Value = "1"
$ws.Cells.Item(35, 10).Value = "0"
$ws.Cells.Item(35, 11).Value = "1"
$ws.Cells.Item(35, 12).Value = "0"
$ws.Cells.Item(35, 13).Value = "0"
$ws.Cells.Item(35, 14).Value = "1"
$ws.Cells.Item(35, 15).Value = "0"

$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(36, 3).Value = "ranking"
$ws.Cells.Item(36, 4).Value = "tensorflow"
$ws.Cells.Item(36, 5).Value = "12/03/2018"
$ws.Cells.Item(36, 6).Value = "0"
$ws.Cells.Item(36, 7).Value = "0"
$ws.Cells.Item(36, 8).Value = "0"
$ws.Cells.Item(36, 9).Value = "1"
$ws.Cells.Item(36, 10).Value = "1"
$ws.Cells.Item(36, 11).Value = "1"
$ws.Cells.Item(36, 12).Value = "0"
$ws.Cells.Item(36, 13).Value = "0"
$ws.Cells.Item(36, 14).Value = "1"
$ws.Cells.Item(36, 15).Value = "0"

$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(37, 3).Value = "ranking"
$ws.Cells.Item(37, 4).Value = "tensorflow"
$ws.Cells.Item(37, 5).Value = "12/03/2018"
$ws.Cells.Item(37, 6).Value = "0"
$ws.Cells.Item(37, 7).Value = "0"
$ws.Cells.Item(37, 8).Value = "0"
$ws.Cells.Item(37, 9).Value = "1"
$ws.Cells.Item(37, 10).Value = "1"
$ws.Cells.Item(37, 11).Value = "1"
$ws.Cells.Item(37, 12).Value = "0"
$ws.Cells.Item(37, 13).Value = "0"
$ws.Cells.Item(37, 14).Value = "1"
$ws.Cells.Item(37, 15).Value = "0"

$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(38, 3).Value = "ranking"
$ws.Cells.Item(38, 4).Value = "tensorflow"
$ws.Cells.Item(38, 5).Value = "12/03/2018"
$ws.Cells.Item(38, 6).Value = "0"
$ws.Cells.Item(38, 7).Value = "0"
$ws.Cells.Item(38, 8).Value = "0"
$ws.Cells.Item(38, 9).Value = "1"
$ws.Cells.Item(38, 10).Value = "0"
$ws.Cells.Item(38, 11).Value = "1"
$ws.Cells.Item(38, 12).Value = "0"
$ws.Cells.Item(38, 13).Value = "0"
$ws.Cells.Item(38, 14).Value = "1"
$ws.Cells.Item(38, 15).Value = "0"

$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(39, 3).Value = "ranking"
$ws.Cells.Item(39, 4).Value = "tensorflow"
$ws.Cells.Item(39, 5).Value = "12/03/2018"
$ws.Cells.Item(39, 6).Value = "0"
$ws.Cells.Item(39, 7).Value = "0"
$ws.Cells.Item(39, 8).Value = "0"
$ws.Cells.Item(39, 9).Value = "1"
$ws.Cells.Item(39, 10).Value = "1"
$ws.Cells.Item(39, 11).Value = "1"
$ws.Cells.Item(39, 12).Value = "0"
$ws.Cells.Item(39, 13).Value = "0"
$ws.Cells.Item(39, 14).Value = "1"
$ws.Cells.Item(39, 15).Value = "0"

$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(40, 3).Value = "ranking"
$ws.Cells.Item(40, 4).Value = "tensorflow"
$ws.Cells.Item(40, 5).Value = "12/03/2018"
$ws.Cells.Item(40, 6).Value = "0"
$ws.Cells.Item(40, 7).Value = "0"
$ws.Cells.Item(40, 8).Value = "0"
$ws.Cells.Item(40, 9).Value = "1"
$ws.Cells.Item(40, 10).Value = "1"
$ws.Cells.Item(40, 11).Value = "1"
$ws.Cells.Item(40, 12).Value = "0"
$ws.Cells.Item(40, 13).Value = "0"
$ws.Cells.Item(40, 14).Value = "1"
$ws.Cells.Item(40, 15).Value = "0"

$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(41, 3).Value = "ranking"
$ws.Cells.Item(41, 4).Value = "tensorflow"
$ws.Cells.Item(41, 5).Value = "12/03/2018"
$ws.Cells.Item(41, 6).Value = "0"
$ws.Cells.Item(41, 7).Value = "0"
$ws.Cells.Item(41, 8).Value = "0"
$ws.Cells.Item(41, 9).Value = "1"
$ws.Cells.Item(41, 10).Value = "1"
$ws.Cells.Item(41, 11).Value = "1"
$ws.Cells.Item(41, 12).Value = "0"
$ws.Cells.Item(41, 13).Value = "0"
$ws.Cells.Item(41, 14).Value = "1"
$ws.Cells.Item(41, 15).Value = "0"

$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(42, 3).Value = "ranking"
$ws.Cells.Item(42, 4).Value = "tensorflow"
$ws.Cells.Item(42, 5).Value = "12/03/2018"
$ws.Cells.Item(42, 6).Value = "0"
$ws.Cells.Item(42, 7).Value = "0"
$ws.Cells.Item(42, 8).Value = "0"
$ws.Cells.Item(42, 9).Value = "1"
$ws.Cells.Item(42, 10).Value = "1"
$ws.Cells.Item(42, 11).Value = "1"
$ws.Cells.Item(42, 12).Value = "0"
$ws.Cells.Item(42, 13).Value = "0"
$ws.Cells.Item(42, 14).Value = "1"
$ws.Cells.Item(42, 15).Value = "0"

$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(43, 3).Value = "ranking"
$ws.Cells.Item(43, 4).Value = "tensorflow"
$ws.Cells.Item(43, 5).Value = "12/03/2018"
$ws.Cells.Item(43, 6).Value = "0"
$ws.Cells.Item(43, 7).Value = "0"
$ws.Cells.Item(43, 8).Value = "0"
$ws.Cells.Item(43, 9).Value = "1"
$ws.Cells.Item(43, 10).Value = "1"
$ws.Cells.Item(43, 11).Value = "1"
$ws.Cells.Item(43, 12).Value = "0"
$ws.Cells.Item(43, 13).Value = "0"
$ws.Cells.Item(43, 14).Value = "1"
$ws.Cells.Item(43, 15).Value = "0"

$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(44, 3).Value = "ranking"
$ws.Cells.Item(44, 4).Value = "tensorflow"
$ws.Cells.Item(44, 5).Value = "12/03/2018"
$ws.Cells.Item(44, 6).Value = "0"
$ws.Cells.Item(44, 7).Value = "0"
$ws.Cells.Item(44, 8).Value = "0"
$ws.Cells.Item(44, 9).Value = "1"
$ws.Cells.Item(44, 10).Value = "1"
$ws.Cells.Item(44, 11).Value = "1"
$ws.Cells.Item(44, 12).Value = "0"
$ws.Cells.Item(44, 13).Value = "0"
$ws.Cells.Item(44, 14).Value = "1"
$ws.Cells.Item(44, 15).Value = "0"

$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(45, 3).Value = "ranking"
$ws.Cells.Item(45, 4).Value = "tensorflow"
$ws.Cells.Item(45, 5).Value = "12/03/2018"
$ws.Cells.Item(45, 6).Value = "0"
$ws.Cells.Item(45, 7).Value = "0"
$ws.Cells.Item(45, 8).Value = "0"
$ws.Cells.Item(45, 9).Value = "1"
$ws.Cells.Item(45, 10).Value = "1"
$ws.Cells.Item(45, 11).Value = "1"
$ws.Cells.Item(45, 12).Value = "0"
$ws.Cells.Item(45, 13).Value = "0"
$ws.Cells.Item(45, 14).Value = "1"
$ws.Cells.Item(45, 15).Value = "0"

$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(46, 3).Value = "ranking"
$ws.Cells.Item(46, 4).Value = "tensorflow"
$ws.Cells.Item(46, 5).Value = "12/03/2018"
$ws.Cells.Item(46, 6).Value = "0"
$ws.Cells.Item(46, 7).Value = "0"
$ws.Cells.Item(46, 8).Value = "0"
$ws.Cells.Item(46, 9).Value = "1"
$ws.Cells.Item(46, 10).Value = "0"
$ws.Cells.Item(46, 11).Value = "1"
$ws.Cells.Item(46, 12).Value = "0"
$ws.Cells.Item(46, 13).Value = "0"
$ws.Cells.Item(46, 14).Value = "1"
$ws.Cells.Item(46, 15).Value = "0"

$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(47, 3).Value = "ranking"
$ws.Cells.Item(47, 4).Value = "tensorflow"
$ws.Cells.Item(47, 5).Value = "12/03/2018"
$ws.Cells.Item(47, 6).Value = "0"
$ws.Cells.Item(47, 7).Value = "0"
$ws.Cells.Item(47, 8).Value = "0"
$ws.Cells.Item(47, 9).Value = "1"
$ws.Cells.Item(47, 10).Value = "1"
$ws.Cells.Item(47, 11).Value = "1"
$ws.Cells.Item(47, 12).Value = "0"
$ws.Cells.Item(47, 13).Value = "0"
$ws.Cells.Item(47, 14).Value = "1"
$ws.Cells.Item(47, 15).Value = "0"

$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(48, 3).Value = "ranking"
$ws.Cells.Item(48, 4).Value = "tensorflow"
$ws.Cells.Item(48, 5).Value = "12/03/2018"
$ws.Cells.Item(48, 6).Value = "0"
$ws.Cells.Item(48, 7).Value = "0"
$ws.Cells.Item(48, 8).Value = "0"
$ws.Cells.Item(48, 9).Value = "1"
$ws.Cells.Item(48, 10).Value = "0"
$ws.Cells.Item(48, 11).Value = "1"
$ws.Cells.Item(48, 12).Value = "0"
$ws.Cells.Item(48, 13).Value = "0"
$ws.Cells.Item(48, 14).Value = "1"
$ws.Cells.Item(48, 15).Value = "0"

$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(49, 3).Value = "ranking"
$ws.Cells.Item(49, 4).Value = "tensorflow"
$ws.Cells.Item(49, 5).Value = "12/03/2018"
$ws.Cells.Item(49, 6).Value = "0"
$ws.Cells.Item(49, 7).Value = "0"
$ws.Cells.Item(49, 8).Value = "0"
$ws.Cells.Item(49, 9).Value = "1"
$ws.Cells.Item(49, 10).Value = "0"
$ws.Cells.Item(49, 11).Value = "1"
$ws.Cells.Item(49, 12).Value = "0"
$ws.Cells.Item(49, 13).Value = "0"
$ws.Cells.Item(49, 14).Value = "1"
$ws.Cells.Item(49, 15).Value = "0"

$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(50, 3).Value = "ranking"
$ws.Cells.Item(50, 4).Value = "tensorflow"
$ws.Cells.Item(50, 5).Value = "12/03/2018"
$ws.Cells.Item(50, 6).Value = "0"
$ws.Cells.Item(50, 7).Value = "0"
$ws.Cells.Item(50, 8).Value = "0"
$ws.Cells.Item(50, 9).Value = "1"
$ws.Cells.Item(50, 10).Value = "0"
$ws.Cells.Item(50, 11).Value = "1"
$ws.Cells.Item(50, 12).Value = "0"
$ws.Cells.Item(50, 13).Value = "0"
$ws.Cells.Item(50, 14).Value = "1"
$ws.Cells.Item(50, 15).Value = "0"

$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(51, 3).Value = "ranking"
$ws.Cells.Item(51, 4).Value = "tensorflow"
$ws.Cells.Item(51, 5).Value = "12/03/2018"
$ws.Cells.Item(51, 6).Value = "0"
$ws.Cells.Item(51, 7).Value = "0"
$ws.Cells.Item(51, 8).Value = "0"
$ws.Cells.Item(51, 9).Value = "1"
$ws.Cells.Item(51, 10).Value = "0"
$ws.Cells.Item(51, 11).Value = "1"
$ws.Cells.Item(51, 12).Value = "0"
$ws.Cells.Item(51, 13).Value = "0"
$ws.Cells.Item(51, 14).Value = "1"
$ws.Cells.Item(51, 15).Value = "0"

$ws.Cells.Item(52, 1).Value = 51
$ws.Cells.Item(52, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(52, 3).Value = "ranking"
$ws.Cells.Item(52, 4).Value = "tensorflow"
$ws.Cells.Item(52, 5).Value = "12/03/2018"
$ws.Cells.Item(52, 6).Value = "0"
$ws.Cells.Item(52, 7).Value = "0"
$ws.Cells.Item(52, 8).Value = "0"
$ws.Cells.Item(52, 9).Value = "1"
$ws.Cells.Item(52, 10).Value = "0"
$ws.Cells.Item(52, 11).Value = "1"
$ws.Cells.Item(52, 12).Value = "0"
$ws.Cells.Item(52, 13).Value = "0"
$ws.Cells.Item(52, 14).Value = "1"
$ws.Cells.Item(52, 15).Value = "0"

$ws.Cells.Item(53, 1).Value = 52
$ws.Cells.Item(53, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(53, 3).Value = "ranking"
$ws.Cells.Item(53, 4).Value = "tensorflow"
$ws.Cells.Item(53, 5).Value = "12/03/2018"
$ws.Cells.Item(53, 6).Value = "0"
$ws.Cells.Item(53, 7).Value = "0"
$ws.Cells.Item(53, 8).Value = "0"
$ws.Cells.Item(53, 9).Value = "1"
$ws.Cells.Item(53, 10).Value = "0"
$ws.Cells.Item(53, 11).Value = "1"
$ws.Cells.Item(53, 12).Value = "0"
$ws.Cells.Item(53, 13).Value = "0"
$ws.Cells.Item(53, 14).Value = "1"
$ws.Cells.Item(53, 15).Value = "0"

$ws.Cells.Item(54, 1).Value = 53
$ws.Cells.Item(54, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(54, 3).Value = "ranking"
$ws.Cells.Item(54, 4).Value = "tensorflow"
$ws.Cells.Item(54, 5).Value = "12/03/2018"
$ws.Cells.Item(54, 6).Value = "0"
$ws.Cells.Item(54, 7).Value = "0"
$ws.Cells.Item(54, 8).Value = "0"
$ws.Cells.Item(54, 9).Value = "1"
$ws.Cells.Item(54, 10).Value = "0"
$ws.Cells.Item(54, 11).Value = "1"
$ws.Cells.Item(54, 12).Value = "0"
$ws.Cells.Item(54, 13).Value = "0"
$ws.Cells.Item(54, 14).Value = "1"
$ws.Cells.Item(54, 15).Value = "0"

$ws.Cells.Item(55, 1).Value = 54
$ws.Cells.Item(55, 2).Value = "https://github.com/tensorflow/ranking"
$ws.Cells.Item(55, 3).Value = "ranking"
$ws.Cells.Item(55, 4).Value = "tensorflow"
$ws.Cells.Item(55, 5).Value = "12/03/2018"
$ws.Cells.Item(55, 6).Value = "0"
$ws.Cells.Item(55, 7).Value = "0"
$ws.Cells.Item(55, 8).Value = "0"
$ws.Cells.Item(55, 9).Value = "1"
$ws.Cells.Item(55, 10).Value = "0"
$ws.Cells.Item(55, 11).Value = "1"
$ws.Cells.Item(55, 12).Value = "0"
$ws.Cells.Item(55, 13).Value = "0"
$ws.Cells.Item(55, 14).Value = "1"
$ws.Cells.Item(55, 15).Value = "0"
